$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.117.77"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "1.891.19"
$ws.Range("E3").Value = "  +1.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5150"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3745"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07221"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9067"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07638"
$ws.Range("D12").ClearFormats()

$ws.Range("D13").Value = "1.896.19"
$ws.Range("E13").Value = "  +1.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.70"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.275"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008500"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.47"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").Value = "27.149.04"
$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.081"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.14%  "

$ws.Range("D22").Value = "2.133.00"
$ws.Range("E22").Value = "  +2.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.57"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.421"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.37"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.792"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.209"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.62"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.967"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.863"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09202"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05085"
$ws.Range("D33").ClearFormats()

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.237"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7698"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.974"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.298"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.612"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5618"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02001"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("E41").Value = "  -0.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.663"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.969"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "117.72"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("E45").Value = "  +3.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4815"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9993"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.594"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.29%  "
